$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 58) at the bottom of the existing table.
$ws.Range("A58").Value = "CompaNanny"
$ws.Range("B58").Value = "CompaNanny Statenkwartier BSO"
$ws.Range("C58").Value = "VGO"

# Column D holds the report date as plain text in this sheet (e.g. the
# preceding rows store "2020-08-13", "2019-11-07", "2024-08-29" as text).
# Force the cell to text first so Excel doesn't auto-convert the
# date-looking string into a real date serial number, then restore the
# default "Normal" cell style so no extra formatting is left behind.
$ws.Range("D58").NumberFormat = "@"
$ws.Range("D58").Value = "2023-04-28"
$ws.Range("D58").Style = "Normal"

$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
